# Update cryptocurrency price/volume data per Sat Nov  9 22:58:01 UTC 2024 GitHub Actions run.
# Rows 8/9 (XRP/Dogecoin) and 43/44 (RenderToken/USDe) are swapped in rank order,
# and D-column price cells are forced to remain plain text (matching the source's
# inlineStr representation, e.g. "76.506.04", "163.80", "1.00") instead of being
# auto-converted to numeric values by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.506.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.126.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.72%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "624.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.68%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.214"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.467"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.09%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.684.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000200"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.435.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.091.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +21.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.46%  "
$ws.Range("E22").Value = "  +5.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.260.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "72.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.08%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  +5.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "514.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.67%  "
$ws.Range("E34").Value = "  +6.80%  "
$ws.Range("E35").Value = "  +19.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "21.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "198.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.66%  "
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.382"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("E42").Value = "  -6.30%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.79%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.804"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +20.50%  "
$ws.Range("E46").Value = "  +7.51%  "
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "41.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.612"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("E51").Value = "  +1.75%  "
